# Update the "provincias Spain" workbook (sheet "Ciudades"):
# 1. Update the "last updated" timestamp text in A1 from 09:50 to 10:20
# 2. Update a few "Muertes" (deaths) figures in column E

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# 1. Update timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 10:20"

# 2. Update death counts for A Coruña (row 16), Pontevedra (row 18), Ourense (row 44)
$ws.Range("E16").Value = 55
$ws.Range("E18").Value = 26
$ws.Range("E44").Value = 16
